# B6-PowerPoint.pptx — re-theme + table style update
#
# The commit swaps the "Office Theme" and "Integral" theme definitions
# between ppt/theme/theme1.xml and ppt/theme/theme2.xml (the deck's live
# design becomes the default "Office Theme" colors instead of the
# "Integral" / "Red Violet" palette), and re-points the three plain
# (Google-Slides-style) tables on the deck to the new table-style GUID.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style id swap for the three tables still using the old,
#    locally-defined "Table_0" style.
# ---------------------------------------------------------------------
$oldStyleId = "{AF2EDBAE-FC42-4ACE-872C-554DDF33EB32}"
$newStyleId = "{DA01D3B7-44F5-4565-9D28-AC53FCADEA11}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            if ($shape.Table.Style -eq $oldStyleId) {
                $shape.Table.ApplyStyle($newStyleId)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Re-theme the deck: replace the live "Integral" / "Red Violet"
#    color scheme with the standard "Office Theme" / "Office" colors.
# ---------------------------------------------------------------------
function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Colors(1..12) == dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le 12; $k++) {
    $themeColors.Colors($k).RGB = HexToComRgb($officeThemeColors[$k - 1])
}
